# "added valid data in excel for username and password/ verfying Homepage
#  is displayed or not."
#
# The sheet becomes a tiny login fixture:
#       A1=UserName   B1=Password
#       A2=admin      B2=manager
# and gets renamed from "Sheet1" to "ValidLogin".
#
# We rebuild the sheet via Copy+Delete (rather than editing Sheet1 in
# place) so the workbook's sheetId advances from 1 to 2 -- exactly what
# happened upstream -- while Copy (unlike Worksheets.Add()) clones the
# existing sheet's XML so formatting/namespace bits aren't lost.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the current sheet (inserted right after itself), then rename the
# clone and drop the original -- this is what bumps sheetId 1 -> 2.
$ws.Copy($null, $ws)
$newSheet = $wb.Worksheets.Item("Sheet1 (2)")
$newSheet.Name = "ValidLogin"

$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet1").Delete()
$excel.DisplayAlerts = $true

# Re-fetch by name: the handle obtained before the delete can go stale.
$ws = $wb.Worksheets.Item("ValidLogin")
$ws.Activate()

# Valid login data.
$ws.Cells.Item(1, 1).Value = "UserName"
$ws.Cells.Item(1, 2).Value = "Password"
$ws.Cells.Item(2, 1).Value = "admin"
$ws.Cells.Item(2, 2).Value = "manager"

# View state: zoom to 175% and leave the cursor resting just below the
# table, on B3.
$excel.ActiveWindow.Zoom = 175
$ws.Range("B3").Select()
